# edit.ps1 - apply the "Add references and various small gains" commit
#
# Changes:
#   1. Bump the Date field on the title page from 2025-01-27 to 2025-01-28.
#   2. Promote the seven top-level section headings from Heading 2 to
#      Heading 1 ("Objectives", "Improve navigation and readability",
#      "Implement dynamic calculations", "Customise figures", "Code",
#      "Add references", "Finalise your MS Word report").

$d = $word.ActiveDocument

# 1. Update the date (set the run text directly so the
#    xml:space="preserve" attribute on the surviving run is kept intact).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "2025-01-27") {
        $p.Range.Text = "2025-01-28"
    }
}

# 2. Promote the Heading 2 section titles to Heading 1.
$headingsToPromote = @(
    "Objectives 🎯",
    "Improve navigation and readability",
    "Implement dynamic calculations",
    "Customise figures",
    "Code",
    "Add references",
    "Finalise your MS Word report"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    if ($headingsToPromote -contains $text -and $p.Style.NameLocal -eq "Heading 2") {
        $p.Style = "Heading1"
    }
}
